# Update countries & provincias Spain
# Applies the data refresh represented in the diff:
#  - footer timestamp text update
#  - re-sorted rows (country labels move, carrying the rest of the row's
#    figures, because the underlying case counts changed and the sheet is
#    kept sorted by total cases)
#  - updated numeric figures for several countries

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer / "last updated" text -------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 09:35"

# --- Rows 28-30: Singapur overtakes Catar & Bielorrusia ----------------
$ws.Cells.Item(28, 1).Value = "Singapur"
$ws.Cells.Item(28, 2).Value = 25346
$ws.Cells.Item(28, 3).Value = 675
$ws.Cells.Item(28, 4).Value = 3851
$ws.Cells.Item(28, 5).Value = 21474
$ws.Cells.Item(28, 6).Value = 20
$ws.Cells.Item(28, 7).Value = 0
$ws.Cells.Item(28, 8).Value = 21

$ws.Cells.Item(29, 1).Value = "Catar"
$ws.Cells.Item(29, 2).Value = 25149
$ws.Cells.Item(29, 3).Value = 0
$ws.Cells.Item(29, 4).Value = 3019
$ws.Cells.Item(29, 5).Value = 22116
$ws.Cells.Item(29, 6).Value = 72
$ws.Cells.Item(29, 7).Value = 0
$ws.Cells.Item(29, 8).Value = 14

$ws.Cells.Item(30, 1).Value = "Bielorrusia"
$ws.Cells.Item(30, 2).Value = 24873
$ws.Cells.Item(30, 3).Value = 0
$ws.Cells.Item(30, 4).Value = 6974
$ws.Cells.Item(30, 5).Value = 17757
$ws.Cells.Item(30, 6).Value = 92
$ws.Cells.Item(30, 7).Value = 0
$ws.Cells.Item(30, 8).Value = 142

# --- Row 11: Alemania figures refreshed --------------------------------
$ws.Cells.Item(11, 2).Value = 173273
$ws.Cells.Item(11, 3).Value = 102
$ws.Cells.Item(11, 5).Value = 16819
$ws.Cells.Item(11, 7).Value = 16
$ws.Cells.Item(11, 8).Value = 7754

# --- Row 33: Polonia figures refreshed ---------------------------------
$ws.Cells.Item(33, 4).Value = 6410
$ws.Cells.Item(33, 5).Value = 9672

# --- Row 39: Israel figures refreshed -----------------------------------
$ws.Cells.Item(39, 5).Value = 7086
$ws.Cells.Item(39, 7).Value = 5
$ws.Cells.Item(39, 8).Value = 1007

# --- Row 51: Chequia figures refreshed -----------------------------------
$ws.Cells.Item(51, 6).Value = 43

# --- Rows 193-194: Nueva Caledonia overtakes Belice ----------------------
$ws.Cells.Item(193, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(193, 2).Value = 18
$ws.Cells.Item(193, 3).Value = 0
$ws.Cells.Item(193, 4).Value = 18
$ws.Cells.Item(193, 5).Value = 0
$ws.Cells.Item(193, 6).Value = 0
$ws.Cells.Item(193, 7).Value = 0
$ws.Cells.Item(193, 8).Value = 0

$ws.Cells.Item(194, 1).Value = "Belice"
$ws.Cells.Item(194, 2).Value = 18
$ws.Cells.Item(194, 3).Value = 0
$ws.Cells.Item(194, 4).Value = 16
$ws.Cells.Item(194, 5).Value = 0
$ws.Cells.Item(194, 6).Value = 0
$ws.Cells.Item(194, 7).Value = 0
$ws.Cells.Item(194, 8).Value = 2

# --- Rows 198-199: Dominica overtakes Curazao -----------------------------
$ws.Cells.Item(198, 1).Value = "Dominica"
$ws.Cells.Item(198, 2).Value = 16
$ws.Cells.Item(198, 3).Value = 0
$ws.Cells.Item(198, 4).Value = 15
$ws.Cells.Item(198, 5).Value = 1
$ws.Cells.Item(198, 6).Value = 0
$ws.Cells.Item(198, 7).Value = 0
$ws.Cells.Item(198, 8).Value = 0

$ws.Cells.Item(199, 1).Value = "Curazao"
$ws.Cells.Item(199, 2).Value = 16
$ws.Cells.Item(199, 3).Value = 0
$ws.Cells.Item(199, 4).Value = 14
$ws.Cells.Item(199, 5).Value = 1
$ws.Cells.Item(199, 6).Value = 0
$ws.Cells.Item(199, 7).Value = 0
$ws.Cells.Item(199, 8).Value = 1

# --- Rows 215-216: San Bartolome overtakes Sahara Occidental --------------
$ws.Cells.Item(215, 1).Value = "San Bartolome"
$ws.Cells.Item(215, 2).Value = 6
$ws.Cells.Item(215, 3).Value = 0
$ws.Cells.Item(215, 4).Value = 6
$ws.Cells.Item(215, 5).Value = 0
$ws.Cells.Item(215, 6).Value = 0
$ws.Cells.Item(215, 7).Value = 0
$ws.Cells.Item(215, 8).Value = 0

$ws.Cells.Item(216, 1).Value = "Sahara Occidental"
$ws.Cells.Item(216, 2).Value = 6
$ws.Cells.Item(216, 3).Value = 0
$ws.Cells.Item(216, 4).Value = 6
$ws.Cells.Item(216, 5).Value = 0
$ws.Cells.Item(216, 6).Value = 0
$ws.Cells.Item(216, 7).Value = 0
$ws.Cells.Item(216, 8).Value = 0
